$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.284.56"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "3.416.60"
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.91"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.98"
$ws.Range("E6").Value = "  +1.38%  "

$ws.Range("E7").Value = "  +3.48%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "3.414.83"
$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.96"
$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("D13").Value = "4.012.69"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.36"
$ws.Range("E15").Value = "  -2.91%  "

$ws.Range("D16").Value = "66.335.24"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000172"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "3.422.25"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.91"
$ws.Range("E19").Value = "  -0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.77"
$ws.Range("E20").Value = "  -0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.38"
$ws.Range("E21").Value = "  -3.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.59"
$ws.Range("E22").Value = "  -3.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.75"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000126"
$ws.Range("E25").Value = "  +5.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.535"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.82"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("E28").Value = "  +1.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.78"
$ws.Range("E30").Value = "  -1.38%  "

$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.42"
$ws.Range("E32").Value = "  -3.68%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.02"
$ws.Range("E34").Value = "  -2.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.26"
$ws.Range("E35").Value = "  -4.12%  "

$ws.Range("E36").Value = "  -1.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.57"
$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.874"
$ws.Range("E38").Value = "  -2.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.63"
$ws.Range("E39").Value = "  -6.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.58"
$ws.Range("E41").Value = "  -2.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.45"
$ws.Range("E42").Value = "  -1.51%  "

$ws.Range("D43").Value = "2.700.75"
$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.36"
$ws.Range("E44").Value = "  -1.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0689"
$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.28"
$ws.Range("E46").Value = "  +2.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.94"
$ws.Range("E47").Value = "  -1.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "333.74"
$ws.Range("E48").Value = "  +8.00%  "

$ws.Range("E49").Value = "  -2.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.104"
$ws.Range("E50").Value = "  +2.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.03"
$ws.Range("E51").Value = "  +4.23%  "
